$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 643.1111
$ws.Range("I33").Value = 155.73334
$ws.Range("K33").Value = 155.73334
$ws.Range("M33").Value = 73.26666
$ws.Range("H62").Value = 2494.3333
$ws.Range("I62").Value = 1934
$ws.Range("K62").Value = 1934
$ws.Range("M62").Value = -1310
$ws.Range("H65").Value = 2494.3333
$ws.Range("I65").Value = 1934
$ws.Range("K65").Value = 9670
$ws.Range("M65").Value = -6550
$ws.Range("H98").Value = 436.8125
$ws.Range("I98").Value = 422.6154
$ws.Range("K98").Value = 422.6154
$ws.Range("M98").Value = 1075.3846
$ws.Range("H107").Value = 484.29413
$ws.Range("I107").Value = 544.4286
$ws.Range("J107").Value = 203.66667
$ws.Range("K107").Value = 544.4286
$ws.Range("L107").Value = 203.66667
$ws.Range("M107").Value = 1375.5714
$ws.Range("N107").Value = -4043.66667
$ws.Range("H116").Value = 3502.6667
$ws.Range("I116").Value = 1475
$ws.Range("K116").Value = 1475
$ws.Range("M116").Value = 1967
$ws.Range("H122").Value = 436.8125
$ws.Range("I122").Value = 422.6154
$ws.Range("K122").Value = 1267.8462
$ws.Range("M122").Value = 1182.1538
$ws.Range("H128").Value = 36328.57
$ws.Range("I128").Value = 31000
$ws.Range("J128").Value = 37216.668
$ws.Range("K128").Value = 31000
$ws.Range("L128").Value = 37216.668
$ws.Range("M128").Value = -26020
$ws.Range("N128").Value = -47176.668
$ws.Range("H132").Value = 5958159.5
$ws.Range("I132").Value = 6416172
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 19248516
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -19245986
$ws.Range("N132").Value = -17059.0001
$ws.Range("H140").Value = 68000
$ws.Range("J140").Value = 68000
$ws.Range("L140").Value = 68000
$ws.Range("N140").Value = -78360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H22").Value = 441
$ws.Range("I22").Value = 441
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 441
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -142
$ws.Range("N22").ClearContents()
$ws.Range("H32").Value = 6229.7124
$ws.Range("I32").Value = 4910.7856
$ws.Range("K32").Value = 4910.7856
$ws.Range("M32").Value = -4623.7856
$ws.Range("H101").Value = 40245
$ws.Range("J101").Value = 40245
$ws.Range("L101").Value = 40245
$ws.Range("N101").Value = -46735
$ws.Range("H112").Value = 9795.666999999999
$ws.Range("J112").Value = 9795.666999999999
$ws.Range("L112").Value = 9795.666999999999
$ws.Range("N112").Value = -12749.667
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H124").Value = 25857.25
$ws.Range("J124").Value = 25857.25
$ws.Range("L124").Value = 25857.25
$ws.Range("N124").Value = -35677.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1804.409
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 1805.1052
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 1805.1052
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -4801.1052
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H112").Value = 32000
$ws.Range("J112").Value = 32000
$ws.Range("L112").Value = 32000
$ws.Range("N112").Value = -34954

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3433.0667
$ws.Range("J86").Value = 3699.4285
$ws.Range("L86").Value = 3699.4285
$ws.Range("N86").Value = -5945.4285
$ws.Range("H89").Value = 3433.0667
$ws.Range("J89").Value = 3699.4285
$ws.Range("L89").Value = 18497.1425
$ws.Range("N89").Value = -29729.1425
$ws.Range("H92").Value = 41326.668
$ws.Range("J92").Value = 41326.668
$ws.Range("L92").Value = 41326.668
$ws.Range("N92").Value = -46318.668
$ws.Range("H122").Value = 600
$ws.Range("I122").Value = 600
$ws.Range("K122").Value = 1800
$ws.Range("M122").Value = 650
$ws.Range("H141").Value = 66003.8
$ws.Range("J141").Value = 45004.75
$ws.Range("L141").Value = 45004.75
$ws.Range("N141").Value = -55364.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1300.25
$ws.Range("I2").Value = 89.25
$ws.Range("J2").Value = 2511.25
$ws.Range("K2").Value = 535.5
$ws.Range("L2").Value = 15067.5
$ws.Range("M2").Value = -422.5
$ws.Range("N2").Value = -15293.5
$ws.Range("H10").Value = 236
$ws.Range("J10").Value = 700
$ws.Range("L10").Value = 2100
$ws.Range("N10").Value = -2378
$ws.Range("H23").Value = 848.25
$ws.Range("I23").Value = 600
$ws.Range("J23").Value = 870.8182
$ws.Range("K23").Value = 1800
$ws.Range("L23").Value = 2612.4546
$ws.Range("M23").Value = -1565
$ws.Range("N23").Value = -3082.4546
$ws.Range("H38").Value = 40.125
$ws.Range("I38").Value = 25.714285
$ws.Range("J38").Value = 51.333332
$ws.Range("K38").Value = 77.142855
$ws.Range("L38").Value = 153.999996
$ws.Range("M38").Value = 269.857145
$ws.Range("N38").Value = -847.999996
$ws.Range("H109").Value = 3123.0417
$ws.Range("I109").Value = 974.875
$ws.Range("J109").Value = 4197.125
$ws.Range("K109").Value = 2924.625
$ws.Range("L109").Value = 12591.375
$ws.Range("M109").Value = -1884.625
$ws.Range("N109").Value = -14671.375
$ws.Range("H112").Value = 92925.45
$ws.Range("J112").Value = 2497.5
$ws.Range("L112").Value = 7492.5
$ws.Range("N112").Value = -9708.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 42666.668
$ws.Range("J111").Value = 42666.668
$ws.Range("L111").Value = 42666.668
$ws.Range("N111").Value = -48800.668
$ws.Range("H132").Value = 1954.4651
$ws.Range("I132").Value = 1745.7715
$ws.Range("J132").Value = 2867.5
$ws.Range("K132").Value = 5237.3145
$ws.Range("L132").Value = 8602.5
$ws.Range("M132").Value = -2707.3145
$ws.Range("N132").Value = -13662.5
$ws.Range("H140").Value = 143894.28
$ws.Range("J140").Value = 143894.28
$ws.Range("L140").Value = 143894.28
$ws.Range("N140").Value = -154254.28

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2002
$ws.Range("I3").Value = 2002
$ws.Range("K3").Value = 2002
$ws.Range("M3").Value = -1890
$ws.Range("H15").Value = 2002
$ws.Range("I15").Value = 2002
$ws.Range("K15").Value = 2002
$ws.Range("M15").Value = -1832
$ws.Range("H47").Value = 15980
$ws.Range("J47").Value = 15980
$ws.Range("L47").Value = 15980
$ws.Range("N47").Value = -16960
$ws.Range("H52").Value = 15980
$ws.Range("J52").Value = 15980
$ws.Range("L52").Value = 15980
$ws.Range("N52").Value = -16446
$ws.Range("H75").Value = 30173
$ws.Range("J75").Value = 30173
$ws.Range("L75").Value = 30173
$ws.Range("N75").Value = -32045
$ws.Range("H78").Value = 30173
$ws.Range("J78").Value = 30173
$ws.Range("L78").Value = 90519
$ws.Range("N78").Value = -99879

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 894.6667
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H70").Value = 10500
$ws.Range("J70").Value = 10500
$ws.Range("L70").Value = 10500
$ws.Range("N70").Value = -11130
$ws.Range("H73").Value = 10500
$ws.Range("J73").Value = 10500
$ws.Range("L73").Value = 10500
$ws.Range("N73").Value = -12684
$ws.Range("H101").Value = 16300
$ws.Range("J101").Value = 16300
$ws.Range("L101").Value = 16300
$ws.Range("N101").Value = -22790
$ws.Range("H136").Value = 1379.5
$ws.Range("I136").Value = 505.74194
$ws.Range("J136").Value = 6796.8
$ws.Range("K136").Value = 1517.22582
$ws.Range("L136").Value = 20390.4
$ws.Range("M136").Value = 1032.77418
$ws.Range("N136").Value = -25490.4
